# Gaussian Quadrature export update for CopperA-HW10.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet/tab (sheet name in workbook.xml)
$ws.Name = "CopperA"

# 2. Tiny re-computed-precision corrections on existing rows (row 13 & 15)
$ws.Range("D13").Value = 0.9958155504106795
$ws.Range("H13").Value = 0.9958155504106795
$ws.Range("J13").Value = 0.9984010437941702
$ws.Range("L13").Value = 0.9897808431246279
$ws.Range("O13").Value = 0.9899436219528295
$ws.Range("P13").Value = 0.9913326706124236

$ws.Range("M15").Value = 0.9042591566025288
$ws.Range("O15").Value = 0.9678583101497081

# 3. Append new row 16 (HexGrid-60degTilt5degRes / index 14)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.560627062854345
$ws.Range("D16").Value = 1.039049773750829
$ws.Range("E16").Value = 1.622924357132956
$ws.Range("F16").Value = 0.7297912232260775
$ws.Range("G16").Value = 1.560627062854345
$ws.Range("H16").Value = 1.039049773750829
$ws.Range("I16").Value = 1.015516809015217
$ws.Range("J16").Value = 1.111178102939896
$ws.Range("K16").Value = 0.8657733098031337
$ws.Range("L16").Value = 0.7045161732785121
$ws.Range("M16").Value = 1.560627062854345
$ws.Range("N16").Value = 1.330987065441892
$ws.Range("O16").Value = 1.238098104241052
$ws.Range("P16").Value = 1.081172101500121

# Match the bold/centered/bordered style used by column A row labels (copy from A15)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
